# Natmi following Dr Hou advice
#
# Re-run of the NATMI ligand-receptor summary for Mfng -> Notch2 after adding a
# third sending-cluster group ("sCs", per Dr Hou). This refreshes the computed
# statistics for the existing Sending-cluster x Target-cluster combinations (now
# based on 3 cell groups instead of 1) and appends the 4 new rows produced by the
# "sCs" sending cluster against each target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One array per output row: row number followed by the 20 column values (A:T).
$tableRows = @(
    ,@(2, "ECs", "Mfng", "Notch2", "ECs", 3, 1, 14.416292, 43.248876, 0.8004770782290026, 0.8004770782290026, 3, 1, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 43.15001387489467, 388.350124874052, 0.02789171034762377, 0.02789171034762378)
    ,@(3, "ECs", "Mfng", "Notch2", "FAPs", 3, 1, 14.416292, 43.248876, 0.8004770782290026, 0.8004770782290026, 3, 1, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 461.253039608504, 4151.277356476536, 0.2981490623623312, 0.2981490623623312)
    ,@(4, "ECs", "Mfng", "Notch2", "M2", 3, 1, 14.416292, 43.248876, 0.8004770782290026, 0.8004770782290026, 3, 1, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 494.3244939995706, 4448.920445996136, 0.319526098979855, 0.319526098979855)
    ,@(5, "ECs", "Mfng", "Notch2", "sCs", 3, 1, 14.416292, 43.248876, 0.8004770782290026, 0.8004770782290026, 3, 1, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 239.6546313660693, 2156.891682294623, 0.1549102065391925, 0.1549102065391925)
    ,@(6, "M2", "Mfng", "Notch2", "ECs", 3, 1, 3.178646666666667, 9.53594, 0.1764971045575167, 0.1764971045575167, 3, 1, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 9.514141900708891, 85.62727710638, 0.006149840203299606, 0.006149840203299606)
    ,@(7, "M2", "Mfng", "Notch2", "FAPs", 3, 1, 3.178646666666667, 9.53594, 0.1764971045575167, 0.1764971045575167, 3, 1, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 101.7016329054267, 915.3146961488401, 0.06573885457146793, 0.06573885457146791)
    ,@(8, "M2", "Mfng", "Notch2", "M2", 3, 1, 3.178646666666667, 9.53594, 0.1764971045575167, 0.1764971045575167, 3, 1, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 108.9935543136489, 980.94198882284, 0.07045227506735571, 0.07045227506735571)
    ,@(9, "M2", "Mfng", "Notch2", "sCs", 3, 1, 3.178646666666667, 9.53594, 0.1764971045575167, 0.1764971045575167, 3, 1, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 52.84142379628445, 475.57281416656, 0.03415613471539348, 0.03415613471539348)
    ,@(10, "sCs", "Mfng", "Notch2", "ECs", 2, 0.6666666666666666, 0.4146863333333333, 1.244059, 0.02302581721348076, 0.02302581721348076, 3, 1, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 1.241215219354778, 11.170936974193, 0.0008023083255008635, 0.0008023083255008636)
    ,@(11, "sCs", "Mfng", "Notch2", "FAPs", 2, 0.6666666666666666, 0.4146863333333333, 1.244059, 0.02302581721348076, 0.02302581721348076, 3, 1, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 13.26799788281933, 119.411980945374, 0.008576292812174343, 0.008576292812174344)
    ,@(12, "sCs", "Mfng", "Notch2", "M2", 2, 0.6666666666666666, 0.4146863333333333, 1.244059, 0.02302581721348076, 0.02302581721348076, 3, 1, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 14.21930215436378, 127.973719389274, 0.009191205782337083, 0.009191205782337083)
    ,@(13, "sCs", "Mfng", "Notch2", "sCs", 2, 0.6666666666666666, 0.4146863333333333, 1.244059, 0.02302581721348076, 0.02302581721348076, 3, 1, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 6.893693631312888, 62.043242681816, 0.004456010293468466, 0.004456010293468467)
)

foreach ($row in $tableRows) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}

